$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 31.705883
$ws.Range("I4").Value = 31.705883
$ws.Range("K4").Value = 31.705883
$ws.Range("M4").Value = 82.294117

$ws.Range("H21").Value = 500
$ws.Range("I21").Value = 500
$ws.Range("K21").Value = 500
$ws.Range("M21").Value = -32

$ws.Range("H23").Value = 500
$ws.Range("I23").Value = 500
$ws.Range("K23").Value = 500
$ws.Range("M23").Value = -266

$ws.Range("H92").Value = 3921.1052
$ws.Range("I92").Value = 4423.154
$ws.Range("K92").Value = 4423.154
$ws.Range("M92").Value = -3175.154

$ws.Range("H98").Value = 8438.583000000001
$ws.Range("I98").Value = 8572.182000000001
$ws.Range("K98").Value = 8572.182000000001
$ws.Range("M98").Value = -7074.182000000001

$ws.Range("H112").Value = 5618.2915
$ws.Range("J112").Value = 3463.2778
$ws.Range("L112").Value = 10389.8334
$ws.Range("N112").Value = -12605.8334

$ws.Range("H122").Value = 8438.583000000001
$ws.Range("I122").Value = 8572.182000000001
$ws.Range("K122").Value = 25716.546
$ws.Range("M122").Value = -23266.546

$ws.Range("H129").Value = 54274.383
$ws.Range("I129").Value = 54274.383
$ws.Range("K129").Value = 162823.149
$ws.Range("M129").Value = -157823.149

$ws.Range("H132").Value = 5519.2656
$ws.Range("I132").Value = 4180.75
$ws.Range("K132").Value = 12542.25
$ws.Range("M132").Value = -10012.25

$ws.Range("H133").Value = 79982
$ws.Range("J133").Value = 79982
$ws.Range("L133").Value = 79982
$ws.Range("N133").Value = -90102

$ws.Range("H135").Value = 3968.5715
$ws.Range("I135").Value = 1704.7368
$ws.Range("K135").Value = 15342.6312
$ws.Range("M135").Value = -12807.6312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1254.0385
$ws.Range("I2").Value = 1012.5
$ws.Range("K2").Value = 1012.5
$ws.Range("M2").Value = -899.5

$ws.Range("H21").Value = 3218.5833
$ws.Range("I21").Value = 328.125
$ws.Range("K21").Value = 328.125
$ws.Range("M21").Value = 45.875

$ws.Range("H32").Value = 6837.8667
$ws.Range("I32").Value = 3799.3584
$ws.Range("J32").Value = 14157.909
$ws.Range("K32").Value = 3799.3584
$ws.Range("L32").Value = 14157.909
$ws.Range("M32").Value = -3512.3584
$ws.Range("N32").Value = -14731.909

$ws.Range("H45").Value = 2982.6
$ws.Range("I45").Value = 3357
$ws.Range("K45").Value = 3357
$ws.Range("M45").Value = -2980

$ws.Range("H116").Value = 1254.0385
$ws.Range("I116").Value = 1012.5
$ws.Range("K116").Value = 1012.5
$ws.Range("M116").Value = 1281.5

$ws.Range("H132").Value = 3586475.8
$ws.Range("I132").Value = 4321.3335
$ws.Range("K132").Value = 12964.0005
$ws.Range("M132").Value = -10434.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1254.0385
$ws.Range("I3").Value = 1012.5
$ws.Range("K3").Value = 1012.5
$ws.Range("M3").Value = -898.5

$ws.Range("H99").Value = 11238.667
$ws.Range("I99").Value = 2763.7778
$ws.Range("K99").Value = 2763.7778
$ws.Range("M99").Value = -1265.7778

$ws.Range("H105").Value = 66668230
$ws.Range("I105").Value = 125000720
$ws.Range("J105").Value = 2528.4285
$ws.Range("K105").Value = 125000720
$ws.Range("L105").Value = 2528.4285
$ws.Range("M105").Value = -124998973
$ws.Range("N105").Value = -6022.4285

$ws.Range("H134").Value = 8499.822
$ws.Range("I134").Value = 3975.7778
$ws.Range("J134").Value = 26596
$ws.Range("K134").Value = 11927.3334
$ws.Range("L134").Value = 79788
$ws.Range("M134").Value = -9392.3334
$ws.Range("N134").Value = -84858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1029.8096
$ws.Range("I22").Value = 539.4375
$ws.Range("K22").Value = 539.4375
$ws.Range("M22").Value = -189.4375

$ws.Range("H31").Value = 12286.075
$ws.Range("I31").Value = 5821.1665
$ws.Range("J31").Value = 21983.438
$ws.Range("K31").Value = 5821.1665
$ws.Range("L31").Value = 21983.438
$ws.Range("M31").Value = -5526.1665
$ws.Range("N31").Value = -22573.438

$ws.Range("H34").Value = 12286.075
$ws.Range("I34").Value = 5821.1665
$ws.Range("J34").Value = 21983.438
$ws.Range("K34").Value = 5821.1665
$ws.Range("L34").Value = 21983.438
$ws.Range("M34").Value = -5619.1665
$ws.Range("N34").Value = -22387.438

$ws.Range("H98").Value = 250000
$ws.Range("J98").Value = 250000
$ws.Range("L98").Value = 250000
$ws.Range("N98").Value = -254492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2987
$ws.Range("I63").Value = 2987
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 8961
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -8212

$ws.Range("H64").Value = 6515.375
$ws.Range("I64").Value = 1056
$ws.Range("K64").Value = 3168
$ws.Range("M64").Value = -2898

$ws.Range("H66").Value = 2987
$ws.Range("I66").Value = 2987
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 26883
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -23139

$ws.Range("H67").Value = 6515.375
$ws.Range("I67").Value = 1056
$ws.Range("K67").Value = 3168
$ws.Range("M67").Value = -2232

$ws.Range("H68").Value = 1271.9474
$ws.Range("I68").Value = 1054.3334
$ws.Range("J68").Value = 1372.3846
$ws.Range("K68").Value = 3163.0002
$ws.Range("L68").Value = 4117.1538
$ws.Range("M68").Value = -2352.0002
$ws.Range("N68").Value = -5739.1538

$ws.Range("H71").Value = 1271.9474
$ws.Range("I71").Value = 1054.3334
$ws.Range("J71").Value = 1372.3846
$ws.Range("K71").Value = 9489.000599999999
$ws.Range("L71").Value = 12351.4614
$ws.Range("M71").Value = -5433.000599999999
$ws.Range("N71").Value = -20463.4614

$ws.Range("H87").Value = 9409.700000000001
$ws.Range("J87").Value = 12999.833
$ws.Range("L87").Value = 38999.499
$ws.Range("N87").Value = -41495.499

$ws.Range("H90").Value = 9409.700000000001
$ws.Range("J90").Value = 12999.833
$ws.Range("L90").Value = 116998.497
$ws.Range("N90").Value = -129478.497

$ws.Range("H134").Value = 6951.5107
$ws.Range("I134").Value = 4872.1
$ws.Range("K134").Value = 14616.3
$ws.Range("M134").Value = -9546.300000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 163.5238
$ws.Range("I2").Value = 156.15384
$ws.Range("J2").Value = 175.5
$ws.Range("K2").Value = 156.15384
$ws.Range("L2").Value = 175.5
$ws.Range("M2").Value = -43.15384
$ws.Range("N2").Value = -401.5

$ws.Range("H102").Value = 6760898
$ws.Range("I102").Value = 9656830
$ws.Range("K102").Value = 9656830
$ws.Range("M102").Value = -9655208

$ws.Range("H126").Value = 1804262.4
$ws.Range("I126").Value = 4000619
$ws.Range("J126").Value = 7243.364
$ws.Range("K126").Value = 12001857
$ws.Range("L126").Value = 21730.092
$ws.Range("M126").Value = -11999387
$ws.Range("N126").Value = -26670.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2832.75
$ws.Range("I61").Value = 2259.9
$ws.Range("K61").Value = 2259.9
$ws.Range("M61").Value = -2057.9

$ws.Range("H82").Value = 2542.5098
$ws.Range("I82").Value = 2620.6365
$ws.Range("J82").Value = 2399.2778
$ws.Range("K82").Value = 2620.6365
$ws.Range("L82").Value = 2399.2778
$ws.Range("M82").Value = -2259.6365
$ws.Range("N82").Value = -3121.2778

$ws.Range("H85").Value = 2542.5098
$ws.Range("I85").Value = 2620.6365
$ws.Range("J85").Value = 2399.2778
$ws.Range("K85").Value = 2620.6365
$ws.Range("L85").Value = 2399.2778
$ws.Range("M85").Value = -1372.6365
$ws.Range("N85").Value = -4895.2778

$ws.Range("H113").Value = 2832.75
$ws.Range("I113").Value = 2259.9
$ws.Range("K113").Value = 2259.9
$ws.Range("M113").Value = -89.90000000000009

$ws.Range("H133").Value = 62442
$ws.Range("J133").Value = 63663
$ws.Range("L133").Value = 63663
$ws.Range("N133").Value = -68723

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7149031
$ws.Range("I126").Value = 7036.25
$ws.Range("J126").Value = 50001000
$ws.Range("K126").Value = 21108.75
$ws.Range("L126").Value = 150003000
$ws.Range("M126").Value = -18638.75
$ws.Range("N126").Value = -150007940

$ws.Range("H132").Value = 5134.6904
$ws.Range("I132").Value = 2766.3
$ws.Range("K132").Value = 8298.900000000001
$ws.Range("M132").Value = -5768.900000000001
